# The commit removes the "H 72" data record (row 2) from the sheet,
# which shifts every following row up by one and shrinks the used
# range from A1:F63 down to A1:F62.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Rows.Item(2).Delete()
